$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = "wor2vec (svm)"
$ws.Range("B10").Value = 0.28999999999999998

$ws.Range("A11").Value = "GloVe (logistic regresion)"
$ws.Range("B11").Value = 0.89

$ws.Range("A12").Value = "Glove (random forest"
$ws.Range("B12").Value = 0.88

$ws.Range("B16").Select() | Out-Null
